$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 24
$ws.Range("H2").Value = "REJET VISA"

$ws.Range("E3").Value = 54
$ws.Range("H3").Value = "REJET VISA"
